$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("result")

$ws.Range("B2").Value = [double]"-20.60755986554534"
$ws.Range("C2").Value = [double]"1.089226581658092"
$ws.Range("D2").Value = [double]"75.00302793421463"
$ws.Range("E2").Value = [double]"-0.406779944259445"
$ws.Range("F2").Value = [double]"-0.352054911594391"
$ws.Range("G2").Value = [double]"0.7815408847023604"
$ws.Range("H2").Value = [double]"-0.4500861918092851"
$ws.Range("I2").Value = [double]"0.8988809271980212"
$ws.Range("J2").Value = [double]"4.35525251498856"
$ws.Range("K2").Value = [double]"27"

$ws.Range("B3").Value = [double]"6.489415692342842e-06"
$ws.Range("C3").Value = [double]"0.7035650081437954"
$ws.Range("D3").Value = [double]"1.195271404766907"
$ws.Range("E3").Value = [double]"62.69770586818203"
$ws.Range("F3").Value = [double]"1.653471968311442"
$ws.Range("G3").Value = [double]"-1.727242461161577"
$ws.Range("H3").Value = [double]"0.5345901828952315"
$ws.Range("I3").Value = [double]"-0.6299647063015685"
$ws.Range("J3").Value = [double]"4.355350721035633"
$ws.Range("K3").Value = [double]"98"

$ws.Range("B4").Value = [double]"-1.836450284879223"
$ws.Range("C4").Value = [double]"62.8767700094226"
$ws.Range("D4").Value = [double]"0.1799469150252427"
$ws.Range("E4").Value = [double]"-0.0001415292567473507"
$ws.Range("F4").Value = [double]"-0.6433970046019519"
$ws.Range("G4").Value = [double]"-0.7616181374878477"
$ws.Range("H4").Value = [double]"0.8604019339048006"
$ws.Range("I4").Value = [double]"1.848922505396178"
$ws.Range("J4").Value = [double]"4.35536275239051"
$ws.Range("K4").Value = [double]"17"

$ws.Range("B5").Value = [double]"0.3725194491593553"
$ws.Range("C5").Value = [double]"8.517556335889296"
$ws.Range("D5").Value = [double]"58.72967226374491"
$ws.Range("E5").Value = [double]"-0.0290272403757472"
$ws.Range("F5").Value = [double]"0.7927274845000563"
$ws.Range("G5").Value = [double]"-1.324838884042288"
$ws.Range("H5").Value = [double]"-0.7077263909663822"
$ws.Range("I5").Value = [double]"1.073934080018848"
$ws.Range("J5").Value = [double]"4.355367582335208"
$ws.Range("K5").Value = [double]"77"

$ws.Range("B6").Value = [double]"-0.3630679222453754"
$ws.Range("C6").Value = [double]"1.444344719566855"
$ws.Range("D6").Value = [double]"47.84348870804013"
$ws.Range("E6").Value = [double]"20.08156574257011"
$ws.Range("F6").Value = [double]"0.8703964381228366"
$ws.Range("G6").Value = [double]"0.7040910487935337"
$ws.Range("H6").Value = [double]"-0.3947333488908169"
$ws.Range("I6").Value = [double]"-0.8932278747492548"
$ws.Range("J6").Value = [double]"4.355476733141792"
$ws.Range("K6").Value = [double]"40"

$ws.Range("B7").Value = [double]"91.69643404025862"
$ws.Range("C7").Value = [double]"-0.5686553812968771"
$ws.Range("D7").Value = [double]"1.817496675291703"
$ws.Range("E7").Value = [double]"558.5564966879476"
$ws.Range("F7").Value = [double]"-1.469611940034627"
$ws.Range("G7").Value = [double]"0.9017877235241905"
$ws.Range("H7").Value = [double]"0.7566475795356018"
$ws.Range("I7").Value = [double]"-0.01278018778117751"
$ws.Range("J7").Value = [double]"4.355514536895612"
$ws.Range("K7").Value = [double]"74"

$ws.Range("B8").Value = [double]"-2.22303467957997"
$ws.Range("C8").Value = [double]"0.002041387304143306"
$ws.Range("D8").Value = [double]"2.123366612837303"
$ws.Range("E8").Value = [double]"72.80241291235619"
$ws.Range("F8").Value = [double]"-0.8754408538818552"
$ws.Range("G8").Value = [double]"1.360473970173114"
$ws.Range("H8").Value = [double]"0.4111161267366548"
$ws.Range("I8").Value = [double]"-0.7137636663214746"
$ws.Range("J8").Value = [double]"4.35551500608134"
$ws.Range("K8").Value = [double]"8"

$ws.Range("B9").Value = [double]"-14.04075639617688"
$ws.Range("C9").Value = [double]"-4.785324361805287e-05"
$ws.Range("D9").Value = [double]"8.052089433750204"
$ws.Range("E9").Value = [double]"58.41453540999616"
$ws.Range("F9").Value = [double]"0.3890331918063996"
$ws.Range("G9").Value = [double]"1.980658151673694"
$ws.Range("H9").Value = [double]"0.4982875334585"
$ws.Range("I9").Value = [double]"-0.993170253252289"
$ws.Range("J9").Value = [double]"4.355557487876269"
$ws.Range("K9").Value = [double]"35"

$ws.Range("B10").Value = [double]"118.8566310237852"
$ws.Range("C10").Value = [double]"384.9131034325741"
$ws.Range("D10").Value = [double]"17.40569882423881"
$ws.Range("E10").Value = [double]"-402.9847387035821"
$ws.Range("F10").Value = [double]"-0.5307320830933719"
$ws.Range("G10").Value = [double]"-1.417852201440257"
$ws.Range("H10").Value = [double]"0.2162672493667328"
$ws.Range("I10").Value = [double]"-1.306742012696025"
$ws.Range("J10").Value = [double]"4.355620433024086"
$ws.Range("K10").Value = [double]"26"

$ws.Range("B11").Value = [double]"-2.402943353734941"
$ws.Range("C11").Value = [double]"44.53439659643279"
$ws.Range("D11").Value = [double]"-0.01328063671396947"
$ws.Range("E11").Value = [double]"2.399261227876057"
$ws.Range("F11").Value = [double]"0.7065799905432026"
$ws.Range("G11").Value = [double]"-0.4627780455198596"
$ws.Range("H11").Value = [double]"1.301177413002297"
$ws.Range("I11").Value = [double]"0.7558072736687058"
$ws.Range("J11").Value = [double]"4.355632374199807"
$ws.Range("K11").Value = [double]"18"

